# Update the NSE hydro cluster table with refreshed figures and append a
# new "Кластер 9" row (row 10).
#
# Every value in columns B:G is stored as plain text (counts like "171",
# medians like "0.76", percentages like "76.6%" are literal strings, not
# numbers) - so we can't just assign $cell.Value = "171" directly, since
# Excel auto-converts number/percent-looking input to Number/Percentage.
# Instead we stage each row's text in a scratch area, copy it, and use
# PasteSpecial(xlPasteValues) into the destination, which carries the
# text across without Excel re-interpreting it and without leaving any
# quote-prefix formatting behind. The scratch columns are deleted again
# once all rows are done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$data = @(
    @{ Row = 2;  B = "171"; C = "0.76"; D = "40"; E = "19"; F = "131"; G = "76.6%" }
    @{ Row = 3;  B = "113"; C = "0.69"; D = "21"; E = "3";  F = "92";  G = "81.4%" }
    @{ Row = 4;  B = "139"; C = "0.83"; D = "11"; E = "2";  F = "128"; G = "92.1%" }
    @{ Row = 5;  B = "156"; C = "0.71"; D = "32"; E = "6";  F = "124"; G = "79.5%" }
    @{ Row = 6;  B = "26";  C = "0.76"; D = "2";  E = "0";  F = "24";  G = "92.3%" }
    @{ Row = 7;  B = "60";  C = "0.81"; D = "8";  E = "1";  F = "52";  G = "86.7%" }
    @{ Row = 8;  B = "102"; C = "0.78"; D = "14"; E = "1";  F = "88";  G = "86.3%" }
    @{ Row = 9;  B = "107"; C = "0.41"; D = "62"; E = "25"; F = "45";  G = "42.1%" }
    @{ Row = 10; B = "122"; C = "0.60"; D = "51"; E = "17"; F = "71";  G = "58.2%" }
)

# Row 10 is new - give column A the same style as the existing cluster
# labels (bold, bordered, centered) by copying the formatting from A9
# before filling in its text.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Cells.Item(10, 1).Value = "Кластер 9"

foreach ($entry in $data) {
    $r = $entry.Row

    $stage = $ws.Range("Z" + $r + ":AE" + $r)
    $stage.Cells.Item(1, 1).Value = "'" + $entry.B
    $stage.Cells.Item(1, 2).Value = "'" + $entry.C
    $stage.Cells.Item(1, 3).Value = "'" + $entry.D
    $stage.Cells.Item(1, 4).Value = "'" + $entry.E
    $stage.Cells.Item(1, 5).Value = "'" + $entry.F
    $stage.Cells.Item(1, 6).Value = "'" + $entry.G

    $stage.Copy()
    $ws.Range("B" + $r + ":G" + $r).PasteSpecial($xlPasteValues)
}

$ws.Range("Z1:AE10").EntireColumn.Delete()
